# Applies the "some design and code changes" edit:
#  - Inserts a new column A (shifting Items -> B, Prices -> C)
#  - Adds a header row: B1 = "Items", C1 = "Prices"
#  - Updates the last price from $1.28 USD to $1.24 USD
#  - Fills A2:A5 with solid red (new style referencing the new fill)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data (A:item name, B:price) one column to the right so it
# becomes (B:item name, C:price), freeing up column A.
$ws.Columns("A").Insert()

# New header row for the shifted columns.
$ws.Range("B1").Value = "Items"
$ws.Range("C1").Value = "Prices"

# Price correction for the last row (Dreams & Nightmares Case).
$ws.Range("C5").Value = "$1.24 USD"

# Highlight the now-empty column A cells next to each data row with a solid
# red fill (255 = pure red in the BGR-packed OLE_COLOR Excel uses).
$ws.Range("A2:A5").Interior.Color = 255
